$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157 (shifts existing rows 157..280 down to 158..281)
$ws.Rows(157).Insert()

# Populate the newly inserted row 157 with the new price-record data
$ws.Cells.Item(157, 1).Value = 4
$ws.Cells.Item(157, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(157, 3).Value = "Los Lagos"
$ws.Cells.Item(157, 4).Value = 44729
$ws.Cells.Item(157, 5).Value = 10
$ws.Cells.Item(157, 6).Value = 100112037
$ws.Cells.Item(157, 7).Value = "Cebollín"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 160
$ws.Cells.Item(157, 11).Value = 11000
$ws.Cells.Item(157, 12).Value = 11000
$ws.Cells.Item(157, 13).Value = 11000
$ws.Cells.Item(157, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(157, 15).Value = "Región Metropolitana"
$ws.Cells.Item(157, 16).Value = 306
$ws.Cells.Item(157, 17).Value = 36
$ws.Cells.Item(157, 18).Value = "Hortaliza"
